$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused rows/cells entirely (formatting + content)
$ws.Range("A35:H37").Clear()
$ws.Range("J25:Q37").Clear()

$data = New-Object "object[,]" 34,17
$data[0,0] = 'negative'
$data[0,1] = $null
$data[0,2] = $null
$data[0,3] = $null
$data[0,4] = $null
$data[0,5] = $null
$data[0,6] = $null
$data[0,7] = $null
$data[0,8] = $null
$data[0,9] = 'positive'
$data[0,10] = $null
$data[0,11] = $null
$data[0,12] = $null
$data[0,13] = $null
$data[0,14] = $null
$data[0,15] = $null
$data[0,16] = $null
$data[1,0] = 'name'
$data[1,1] = 'anchor score'
$data[1,2] = 'type occurences'
$data[1,3] = 'total occurences'
$data[1,4] = '+%'
$data[1,5] = '-%'
$data[1,6] = 'both'
$data[1,7] = 'normal'
$data[1,8] = $null
$data[1,9] = 'name'
$data[1,10] = 'anchor score'
$data[1,11] = 'type occurences'
$data[1,12] = 'total occurences'
$data[1,13] = '+%'
$data[1,14] = '-%'
$data[1,15] = 'both'
$data[1,16] = 'normal'
$data[2,0] = 'poorly'
$data[2,1] = 0.9347826086956522
$data[2,2] = 43
$data[2,3] = 43
$data[2,4] = 0
$data[2,5] = 1
$data[2,6] = $false
$data[2,7] = 3
$data[2,8] = $null
$data[2,9] = 'wonderful'
$data[2,10] = 0.8928571428571429
$data[2,11] = 50
$data[2,12] = 50
$data[2,13] = 1
$data[2,14] = 0
$data[2,15] = $false
$data[2,16] = 6
$data[3,0] = 'disappointing'
$data[3,1] = 0.8636363636363636
$data[3,2] = 38
$data[3,3] = 38
$data[3,4] = 0
$data[3,5] = 1
$data[3,6] = $false
$data[3,7] = 6
$data[3,8] = $null
$data[3,9] = 'awesome'
$data[3,10] = 0.8923076923076924
$data[3,11] = 58
$data[3,12] = 58
$data[3,13] = 1
$data[3,14] = 0
$data[3,15] = $false
$data[3,16] = 7
$data[4,0] = 'poor'
$data[4,1] = 0.7605633802816901
$data[4,2] = 54
$data[4,3] = 54
$data[4,4] = 0
$data[4,5] = 1
$data[4,6] = $false
$data[4,7] = 17
$data[4,8] = $null
$data[4,9] = 'favorite'
$data[4,10] = 0.8602150537634409
$data[4,11] = 80
$data[4,12] = 80
$data[4,13] = 1
$data[4,14] = 0
$data[4,15] = $false
$data[4,16] = 13
$data[5,0] = 'however'
$data[5,1] = 0.703125
$data[5,2] = 45
$data[5,3] = 45
$data[5,4] = 0
$data[5,5] = 1
$data[5,6] = $false
$data[5,7] = 19
$data[5,8] = $null
$data[5,9] = 'excellent'
$data[5,10] = 0.765625
$data[5,11] = 49
$data[5,12] = 49
$data[5,13] = 1
$data[5,14] = 0
$data[5,15] = $false
$data[5,16] = 15
$data[6,0] = 'disappointed'
$data[6,1] = 0.6989247311827957
$data[6,2] = 130
$data[6,3] = 130
$data[6,4] = 0
$data[6,5] = 1
$data[6,6] = $false
$data[6,7] = 56
$data[6,8] = $null
$data[6,9] = 'classic'
$data[6,10] = 0.6981132075471698
$data[6,11] = 37
$data[6,12] = 37
$data[6,13] = 1
$data[6,14] = 0
$data[6,15] = $false
$data[6,16] = 16
$data[7,0] = 'junk'
$data[7,1] = 0.6727272727272727
$data[7,2] = 37
$data[7,3] = 37
$data[7,4] = 0
$data[7,5] = 1
$data[7,6] = $false
$data[7,7] = 18
$data[7,8] = $null
$data[7,9] = 'love'
$data[7,10] = 0.5810616929698709
$data[7,11] = 405
$data[7,12] = 405
$data[7,13] = 1
$data[7,14] = 0
$data[7,15] = $false
$data[7,16] = 292
$data[8,0] = 'broke'
$data[8,1] = 0.6601941747572816
$data[8,2] = 136
$data[8,3] = 136
$data[8,4] = 0
$data[8,5] = 1
$data[8,6] = $false
$data[8,7] = 70
$data[8,8] = $null
$data[8,9] = 'loves'
$data[8,10] = 0.524896265560166
$data[8,11] = 253
$data[8,12] = 253
$data[8,13] = 1
$data[8,14] = 0
$data[8,15] = $false
$data[8,16] = 229
$data[9,0] = 'waste'
$data[9,1] = 0.6418918918918919
$data[9,2] = 95
$data[9,3] = 95
$data[9,4] = 0
$data[9,5] = 1
$data[9,6] = $false
$data[9,7] = 53
$data[9,8] = $null
$data[9,9] = 'thank'
$data[9,10] = 0.5217391304347826
$data[9,11] = 36
$data[9,12] = 36
$data[9,13] = 1
$data[9,14] = 0
$data[9,15] = $false
$data[9,16] = 33
$data[10,0] = 'smaller'
$data[10,1] = 0.5798319327731093
$data[10,2] = 69
$data[10,3] = 69
$data[10,4] = 0
$data[10,5] = 1
$data[10,6] = $false
$data[10,7] = 50
$data[10,8] = $null
$data[10,9] = 'great'
$data[10,10] = 0.4950819672131148
$data[10,11] = 604
$data[10,12] = 604
$data[10,13] = 1
$data[10,14] = 0
$data[10,15] = $false
$data[10,16] = 616
$data[11,0] = 'small'
$data[11,1] = 0.5043478260869565
$data[11,2] = 174
$data[11,3] = 174
$data[11,4] = 0
$data[11,5] = 1
$data[11,6] = $false
$data[11,7] = 171
$data[11,8] = $null
$data[11,9] = 'friends'
$data[11,10] = 0.3968253968253968
$data[11,11] = 75
$data[11,12] = 75
$data[11,13] = 1
$data[11,14] = 0
$data[11,15] = $false
$data[11,16] = 114
$data[12,0] = 'broken'
$data[12,1] = 0.4819277108433735
$data[12,2] = 40
$data[12,3] = 40
$data[12,4] = 0
$data[12,5] = 1
$data[12,6] = $false
$data[12,7] = 43
$data[12,8] = $null
$data[12,9] = 'loved'
$data[12,10] = 0.3516819571865443
$data[12,11] = 115
$data[12,12] = 115
$data[12,13] = 1
$data[12,14] = 0
$data[12,15] = $false
$data[12,16] = 212
$data[13,0] = 'plastic'
$data[13,1] = 0.4251968503937008
$data[13,2] = 54
$data[13,3] = 54
$data[13,4] = 0
$data[13,5] = 1
$data[13,6] = $false
$data[13,7] = 73
$data[13,8] = $null
$data[13,9] = 'best'
$data[13,10] = 0.3333333333333333
$data[13,11] = 40
$data[13,12] = 40
$data[13,13] = 1
$data[13,14] = 0
$data[13,15] = $false
$data[13,16] = 80
$data[14,0] = 'apart'
$data[14,1] = 0.4210526315789473
$data[14,2] = 40
$data[14,3] = 40
$data[14,4] = 0
$data[14,5] = 1
$data[14,6] = $false
$data[14,7] = 55
$data[14,8] = $null
$data[14,9] = 'perfect'
$data[14,10] = 0.3192771084337349
$data[14,11] = 53
$data[14,12] = 53
$data[14,13] = 1
$data[14,14] = 0
$data[14,15] = $false
$data[14,16] = 113
$data[15,0] = 'ok'
$data[15,1] = 0.3984375
$data[15,2] = 51
$data[15,3] = 51
$data[15,4] = 0
$data[15,5] = 1
$data[15,6] = $false
$data[15,7] = 77
$data[15,8] = $null
$data[15,9] = 'learn'
$data[15,10] = 0.2421875
$data[15,11] = 31
$data[15,12] = 31
$data[15,13] = 1
$data[15,14] = 0
$data[15,15] = $false
$data[15,16] = 97
$data[16,0] = 'cheap'
$data[16,1] = 0.3933649289099526
$data[16,2] = 83
$data[16,3] = 83
$data[16,4] = 0
$data[16,5] = 1
$data[16,6] = $false
$data[16,7] = 128
$data[16,8] = $null
$data[16,9] = 'happy'
$data[16,10] = 0.2307692307692308
$data[16,11] = 33
$data[16,12] = 33
$data[16,13] = 1
$data[16,14] = 0
$data[16,15] = $false
$data[16,16] = 110
$data[17,0] = 'thought'
$data[17,1] = 0.3267326732673267
$data[17,2] = 66
$data[17,3] = 66
$data[17,4] = 0
$data[17,5] = 1
$data[17,6] = $false
$data[17,7] = 136
$data[17,8] = $null
$data[17,9] = 'enjoy'
$data[17,10] = 0.2258064516129032
$data[17,11] = 42
$data[17,12] = 42
$data[17,13] = 1
$data[17,14] = 0
$data[17,15] = $false
$data[17,16] = 144
$data[18,0] = 'though'
$data[18,1] = 0.282051282051282
$data[18,2] = 33
$data[18,3] = 33
$data[18,4] = 0
$data[18,5] = 1
$data[18,6] = $false
$data[18,7] = 84
$data[18,8] = $null
$data[18,9] = 'christmas'
$data[18,10] = 0.2168674698795181
$data[18,11] = 54
$data[18,12] = 54
$data[18,13] = 1
$data[18,14] = 0
$data[18,15] = $false
$data[18,16] = 195
$data[19,0] = 'size'
$data[19,1] = 0.2422680412371134
$data[19,2] = 47
$data[19,3] = 47
$data[19,4] = 0
$data[19,5] = 1
$data[19,6] = $false
$data[19,7] = 147
$data[19,8] = $null
$data[19,9] = 'fun'
$data[19,10] = 0.1850877192982456
$data[19,11] = 211
$data[19,12] = 212
$data[19,13] = 1
$data[19,14] = 0
$data[19,15] = $true
$data[19,16] = 929
$data[20,0] = 'hard'
$data[20,1] = 0.215
$data[20,2] = 43
$data[20,3] = 43
$data[20,4] = 0
$data[20,5] = 1
$data[20,6] = $false
$data[20,7] = 157
$data[20,8] = $null
$data[20,9] = 'game'
$data[20,10] = 0.1187540558079169
$data[20,11] = 183
$data[20,12] = 183
$data[20,13] = 1
$data[20,14] = 0
$data[20,15] = $false
$data[20,16] = 1358
$data[21,0] = 'money'
$data[21,1] = 0.2120253164556962
$data[21,2] = 67
$data[21,3] = 67
$data[21,4] = 0
$data[21,5] = 1
$data[21,6] = $false
$data[21,7] = 249
$data[21,8] = $null
$data[21,9] = 'family'
$data[21,10] = 0.09749303621169916
$data[21,11] = 35
$data[21,12] = 35
$data[21,13] = 1
$data[21,14] = 0
$data[21,15] = $false
$data[21,16] = 324
$data[22,0] = 'item'
$data[22,1] = 0.1956521739130435
$data[22,2] = 54
$data[22,3] = 54
$data[22,4] = 0
$data[22,5] = 1
$data[22,6] = $false
$data[22,7] = 222
$data[22,8] = $null
$data[22,9] = 'easy'
$data[22,10] = 0.09090909090909091
$data[22,11] = 34
$data[22,12] = 34
$data[22,13] = 1
$data[22,14] = 0
$data[22,15] = $false
$data[22,16] = 340
$data[23,0] = 'work'
$data[23,1] = 0.1809523809523809
$data[23,2] = 57
$data[23,3] = 58
$data[23,4] = 0.02
$data[23,5] = 0.98
$data[23,6] = $true
$data[23,7] = 258
$data[23,8] = $null
$data[23,9] = 'play'
$data[23,10] = 0.05333333333333334
$data[23,11] = 40
$data[23,12] = 42
$data[23,13] = 0.95
$data[23,14] = 0.05000000000000004
$data[23,15] = $true
$data[23,16] = 710
$data[24,0] = 'would'
$data[24,1] = 0.1780415430267062
$data[24,2] = 120
$data[24,3] = 120
$data[24,4] = 0
$data[24,5] = 1
$data[24,6] = $false
$data[24,7] = 554
$data[24,8] = $null
$data[24,9] = $null
$data[24,10] = $null
$data[24,11] = $null
$data[24,12] = $null
$data[24,13] = $null
$data[24,14] = $null
$data[24,15] = $null
$data[24,16] = $null
$data[25,0] = 'price'
$data[25,1] = 0.1695402298850575
$data[25,2] = 59
$data[25,3] = 59
$data[25,4] = 0
$data[25,5] = 1
$data[25,6] = $false
$data[25,7] = 289
$data[25,8] = $null
$data[25,9] = $null
$data[25,10] = $null
$data[25,11] = $null
$data[25,12] = $null
$data[25,13] = $null
$data[25,14] = $null
$data[25,15] = $null
$data[25,16] = $null
$data[26,0] = 'better'
$data[26,1] = 0.1448598130841121
$data[26,2] = 31
$data[26,3] = 31
$data[26,4] = 0
$data[26,5] = 1
$data[26,6] = $false
$data[26,7] = 183
$data[26,8] = $null
$data[26,9] = $null
$data[26,10] = $null
$data[26,11] = $null
$data[26,12] = $null
$data[26,13] = $null
$data[26,14] = $null
$data[26,15] = $null
$data[26,16] = $null
$data[27,0] = 'product'
$data[27,1] = 0.1365638766519824
$data[27,2] = 62
$data[27,3] = 62
$data[27,4] = 0
$data[27,5] = 1
$data[27,6] = $false
$data[27,7] = 392
$data[27,8] = $null
$data[27,9] = $null
$data[27,10] = $null
$data[27,11] = $null
$data[27,12] = $null
$data[27,13] = $null
$data[27,14] = $null
$data[27,15] = $null
$data[27,16] = $null
$data[28,0] = '2'
$data[28,1] = 0.1169811320754717
$data[28,2] = 31
$data[28,3] = 33
$data[28,4] = 0.06
$data[28,5] = 0.94
$data[28,6] = $true
$data[28,7] = 234
$data[28,8] = $null
$data[28,9] = $null
$data[28,10] = $null
$data[28,11] = $null
$data[28,12] = $null
$data[28,13] = $null
$data[28,14] = $null
$data[28,15] = $null
$data[28,16] = $null
$data[29,0] = 'buy'
$data[29,1] = 0.09295774647887324
$data[29,2] = 33
$data[29,3] = 33
$data[29,4] = 0
$data[29,5] = 1
$data[29,6] = $false
$data[29,7] = 322
$data[29,8] = $null
$data[29,9] = $null
$data[29,10] = $null
$data[29,11] = $null
$data[29,12] = $null
$data[29,13] = $null
$data[29,14] = $null
$data[29,15] = $null
$data[29,16] = $null
$data[30,0] = 'little'
$data[30,1] = 0.09213483146067415
$data[30,2] = 41
$data[30,3] = 45
$data[30,4] = 0.09
$data[30,5] = 0.91
$data[30,6] = $true
$data[30,7] = 404
$data[30,8] = $null
$data[30,9] = $null
$data[30,10] = $null
$data[30,11] = $null
$data[30,12] = $null
$data[30,13] = $null
$data[30,14] = $null
$data[30,15] = $null
$data[30,16] = $null
$data[31,0] = 'like'
$data[31,1] = 0.06919275123558484
$data[31,2] = 42
$data[31,3] = 43
$data[31,4] = 0.02
$data[31,5] = 0.98
$data[31,6] = $true
$data[31,7] = 565
$data[31,8] = $null
$data[31,9] = $null
$data[31,10] = $null
$data[31,11] = $null
$data[31,12] = $null
$data[31,13] = $null
$data[31,14] = $null
$data[31,15] = $null
$data[31,16] = $null
$data[32,0] = 'one'
$data[32,1] = 0.05236270753512133
$data[32,2] = 41
$data[32,3] = 52
$data[32,4] = 0.21
$data[32,5] = 0.79
$data[32,6] = $true
$data[32,7] = 742
$data[32,8] = $null
$data[32,9] = $null
$data[32,10] = $null
$data[32,11] = $null
$data[32,12] = $null
$data[32,13] = $null
$data[32,14] = $null
$data[32,15] = $null
$data[32,16] = $null
$data[33,0] = 'toy'
$data[33,1] = 0.0460122699386503
$data[33,2] = 30
$data[33,3] = 33
$data[33,4] = 0.09
$data[33,5] = 0.91
$data[33,6] = $true
$data[33,7] = 622
$data[33,8] = $null
$data[33,9] = $null
$data[33,10] = $null
$data[33,11] = $null
$data[33,12] = $null
$data[33,13] = $null
$data[33,14] = $null
$data[33,15] = $null
$data[33,16] = $null

$ws.Range("A1:Q34").Value = $data

